$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.537.24'
$ws.Range("E2").Value = '  +1.39%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.917.43'
$ws.Range("E3").Value = '  +0.36%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.52'
$ws.Range("E5").Value = '  -2.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9998'
$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4787'
$ws.Range("E7").Value = '  +2.92%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4101'
$ws.Range("E8").Value = '  +0.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.70'
$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("E10").Value = '  +0.50%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.012'
$ws.Range("E11").Value = '  +0.25%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.42'
$ws.Range("E12").Value = '  +2.29%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.900.15'
$ws.Range("E13").Value = '  -0.46%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.953'
$ws.Range("E14").Value = '  -0.05%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.161'
$ws.Range("E15").Value = '  +0.74%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.61'
$ws.Range("E16").Value = '  +0.32%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  -0.04%  '

$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06616'
$ws.Range("E18").Value = '  +0.52%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001030'
$ws.Range("E19").Value = '  -0.32%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.77'
$ws.Range("E20").Value = '  +1.21%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9997'
$ws.Range("E21").Value = '  -0.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '29.550.85'
$ws.Range("E22").Value = '  +1.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.548'
$ws.Range("E23").Value = '  +1.94%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.52'
$ws.Range("E24").Value = '  +1.89%  '

$ws.Range("E25").Value = '  -1.37%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.131.40'
$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.28'
$ws.Range("E27").Value = '  -2.95%  '

$ws.Range("E28").Value = '  +0.71%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.824'
$ws.Range("E29").Value = '  +6.98%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.139'
$ws.Range("E30").Value = '  +1.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.93'
$ws.Range("E31").Value = '  -0.57%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.056'
$ws.Range("E32").Value = '  +6.72%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09544'
$ws.Range("E33").Value = '  +1.31%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.429'
$ws.Range("E34").Value = '  -0.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.571'
$ws.Range("E35").Value = '  -0.62%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.391'
$ws.Range("E36").Value = '  +1.25%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06105'
$ws.Range("E37").Value = '  -0.07%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02257'
$ws.Range("E38").Value = '  +0.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.342'
$ws.Range("E39").Value = '  -0.55%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.176'
$ws.Range("E40").Value = '  +0.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5890'
$ws.Range("E41").Value = '  +1.17%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.536'
$ws.Range("E42").Value = '  +7.20%  '

$ws.Range("E43").Value = '  +1.08%  '

$ws.Range("E44").Value = '  -0.53%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.08003'
$ws.Range("E45").Value = '  +13.12%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.290'
$ws.Range("E46").Value = '  +1.90%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5551'
$ws.Range("E47").Value = '  +0.60%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '12.13'
$ws.Range("E48").Value = '  +0.17%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.934'
$ws.Range("E49").Value = '  +0.49%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '113.40'
$ws.Range("E50").Value = '  +1.53%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.72'
$ws.Range("E51").Value = '  -6.58%  '
